$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from before the "I 1888 gikk han..."
#    paragraph to the end of the "Han og fire andre..." paragraph (right
#    after the new "." run that closes that sentence).
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Close the "Han og fire andre, var de første menneskene på sørpolen"
# sentence with a period.
$p3 = $d.Paragraphs(3)
$insPt = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$insPt.InsertAfter(".")

# Re-seat the _GoBack bookmark right after that new "." (i.e. at the very
# end of the paragraph's text, before its paragraph mark). A collapsed
# Range sitting exactly at "end of paragraph content" confuses
# Bookmarks.Add on this host, so we work around it: insert a throwaway
# run after the insertion point, add the bookmark just before that run,
# then delete the throwaway run again, leaving the bookmark in place.
$p3b = $d.Paragraphs(3)
$tempPt = $d.Range($p3b.Range.End - 1, $p3b.Range.End - 1)
$tempPt.InsertAfter("TEMP_MARK")

$p3c = $d.Paragraphs(3)
$bmPos = $p3c.Range.End - 1 - 9   # 9 == Len("TEMP_MARK")
$bmRng = $d.Range($bmPos, $bmPos)
$bmRng.Bookmarks.Add("_GoBack")

$p3d = $d.Paragraphs(3)
$delRng = $d.Range($p3d.Range.End - 1 - 9, $p3d.Range.End - 1)
$delRng.Delete()

# ---------------------------------------------------------------------------
# 2) Small typo / punctuation fixes.
# ---------------------------------------------------------------------------

# Close out a few sentences that were missing their final period.
$d.Content.Find.Execute("Han var aldri gift", $false, $false, $false, $false, $false, $true, 1, $false, "Han var aldri gift.", 2) | Out-Null

$d.Content.Find.Execute("Peter Markham Scott)", $false, $false, $false, $false, $false, $true, 1, $false, "Peter Markham Scott).", 2) | Out-Null

$d.Content.Find.Execute("ekspedisjonen i 1914-16", $false, $false, $false, $false, $false, $true, 1, $false, "ekspedisjonen i 1914-16.", 2) | Out-Null

$d.Content.Find.Execute("hver eneste person på skipet overlevde", $false, $false, $false, $false, $false, $true, 1, $false, "hver eneste person på skipet overlevde.", 2) | Out-Null

# Shackleton had three children, not one -- fix the name(s) listed.
$d.Content.Find.Execute("Edward Shackleton", $false, $false, $false, $false, $false, $true, 1, $false, "Edward, Raymond and Cecily", 2) | Out-Null

# Remove the erroneous "heldigvis," (the expedition wasn't lucky -- he got
# sick and couldn't go).
$d.Content.Find.Execute("men heldigvis, ble han syk", $false, $false, $false, $false, $false, $true, 1, $false, "men ble han syk", 2) | Out-Null

Write-Host "Done"
